$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A98").Value = 'Create Nationality and Delete'
$ws.Range("B98").Value = 'PASSED'
$ws.Range("C98").Value = 'chrome'
$ws.Range("A99").Value = 'Create Nationality and Delete'
$ws.Range("B99").Value = 'PASSED'
$ws.Range("C99").Value = 'chrome'
$ws.Range("A100").Value = 'Create Nationality and Delete'
$ws.Range("B100").Value = 'PASSED'
$ws.Range("C100").Value = 'chrome'
$ws.Range("A101").Value = 'Create Nationality and Delete'
$ws.Range("B101").Value = 'PASSED'
$ws.Range("C101").Value = 'chrome'
$ws.Range("A102").Value = 'Create Nationality and Delete'
$ws.Range("B102").Value = 'PASSED'
$ws.Range("C102").Value = 'chrome'
$ws.Range("A103").Value = 'Create Nationality and Delete'
$ws.Range("B103").Value = 'PASSED'
$ws.Range("C103").Value = 'chrome'
$ws.Range("A104").Value = 'Create Nationality and Delete'
$ws.Range("B104").Value = 'PASSED'
$ws.Range("C104").Value = 'chrome'
$ws.Range("A105").Value = 'Create Nationality and Delete'
$ws.Range("B105").Value = 'PASSED'
$ws.Range("C105").Value = 'chrome'
$ws.Range("A106").Value = 'Create a Citizenship'
$ws.Range("B106").Value = 'PASSED'
$ws.Range("C106").Value = 'chrome'
$ws.Range("A107").Value = 'Create a Citizenship'
$ws.Range("B107").Value = 'PASSED'
$ws.Range("C107").Value = 'chrome'
$ws.Range("A108").Value = 'Exam Create and Delete'
$ws.Range("B108").Value = 'FAILED'
$ws.Range("C108").Value = 'chrome'
$ws.Range("A109").Value = 'Exam Create and Delete'
$ws.Range("B109").Value = 'FAILED'
$ws.Range("C109").Value = 'chrome'
$ws.Range("A110").Value = 'Create Inventory and Delete'
$ws.Range("B110").Value = 'UNDEFINED'
$ws.Range("A111").Value = 'Create Inventory and Delete'
$ws.Range("B111").Value = 'FAILED'
$ws.Range("C111").Value = 'chrome'
$ws.Range("A112").Value = 'Create Inventory and Delete'
$ws.Range("B112").Value = 'FAILED'
$ws.Range("C112").Value = 'chrome'
$ws.Range("A113").Value = 'Create Inventory and Delete'
$ws.Range("B113").Value = 'FAILED'
$ws.Range("C113").Value = 'chrome'
$ws.Range("A114").Value = 'Create Inventory and Delete'
$ws.Range("B114").Value = 'FAILED'
$ws.Range("C114").Value = 'chrome'
$ws.Range("A115").Value = 'Create Inventory and Delete'
$ws.Range("B115").Value = 'FAILED'
$ws.Range("C115").Value = 'chrome'
$ws.Range("A116").Value = 'Create Inventory and Delete'
$ws.Range("B116").Value = 'FAILED'
$ws.Range("C116").Value = 'chrome'
$ws.Range("A117").Value = 'Create Inventory and Delete'
$ws.Range("B117").Value = 'FAILED'
$ws.Range("C117").Value = 'chrome'
$ws.Range("A118").Value = 'Create Inventory and Delete'
$ws.Range("B118").Value = 'FAILED'
$ws.Range("C118").Value = 'chrome'
$ws.Range("A119").Value = 'Create Inventory and Delete'
$ws.Range("B119").Value = 'PASSED'
$ws.Range("C119").Value = 'chrome'
$ws.Range("A120").Value = 'Create Inventory and Delete'
$ws.Range("B120").Value = 'FAILED'
$ws.Range("C120").Value = 'chrome'
$ws.Range("A121").Value = 'Create Inventory and Delete'
$ws.Range("B121").Value = 'FAILED'
$ws.Range("C121").Value = 'chrome'
$ws.Range("A122").Value = 'Create Inventory and Delete'
$ws.Range("B122").Value = 'FAILED'
$ws.Range("C122").Value = 'chrome'
$ws.Range("A123").Value = 'Create Inventory and Delete'
$ws.Range("B123").Value = 'PASSED'
$ws.Range("C123").Value = 'chrome'
$ws.Range("A124").Value = 'Create Inventory and Delete'
$ws.Range("B124").Value = 'PASSED'
$ws.Range("C124").Value = 'chrome'
$ws.Range("A125").Value = 'Create Inventory and Delete'
$ws.Range("B125").Value = 'PASSED'
$ws.Range("C125").Value = 'chrome'
$ws.Range("A126").Value = 'Create Inventory and Delete'
$ws.Range("B126").Value = 'PASSED'
$ws.Range("C126").Value = 'chrome'
$ws.Range("A127").Value = 'Create Inventory and Delete'
$ws.Range("B127").Value = 'PASSED'
$ws.Range("C127").Value = 'chrome'
